# Append the source reference for the "Quick start" text at the end of the
# document: one blank paragraph followed by a paragraph containing the
# reference URL, both styled like the surrounding body text (Helvetica
# Neue, 16pt, color #191C1F).

$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Last
$endRange = $lastPara.Range.Duplicate()
$endRange.Collapse(0)

# Inserting right after the existing last paragraph's text (before its
# paragraph mark) lets the new runs inherit that paragraph's character
# formatting automatically, so the new paragraphs come out with the same
# rFonts/size/color as the rest of the body text.
$endRange.InsertAfter("`r`rref: https://getbootstrap.com/docs/4.0/getting-started/introduction/")

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
